$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhDFormulationsDataset_2023")

# --- 1. New shared string used by the header (W1): "Newtonian" replaces "Newtonain ".
$ws.Cells.Item(1, 23).Value = "Newtonian"

# --- 2. Copy formatting (styles) from the last existing data row (253) down through
#        the new rows (254-277) so every column keeps its established number format /
#        alignment (A -> s10, B:S -> s7, U/V -> s2, W -> s12, T unstyled).
$ws.Range("A253:W253").Copy()
$ws.Range("A254:W277").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Bulk-write the literal measurement data for columns A (ID) through T (stability bool).
$data = New-Object 'object[,]' 24,20
$data[0,0] = 253
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 12.295247162772799
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 8.2751854684037607
$data[0,13] = 2.5244166959849599
$data[0,14] = 0
$data[0,15] = 0
$data[0,16] = 0
$data[0,17] = 1.0352734146981399
$data[0,18] = 0
$data[0,19] = $true
$data[1,0] = 254
$data[1,1] = 0
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 6.7450675030102101
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 0
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 13.303807174445399
$data[1,13] = 2.9220736289545202
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0
$data[1,17] = 1.42012022984085
$data[1,18] = 0
$data[1,19] = $true
$data[2,0] = 255
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 8.7695469536642907
$data[2,11] = 0
$data[2,12] = 10.544787136032699
$data[2,13] = 1.4440646907389301
$data[2,14] = 0
$data[2,15] = 0
$data[2,16] = 0
$data[2,17] = 1.2210383062330501
$data[2,18] = 0
$data[2,19] = $false
$data[3,0] = 256
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 13.8205367068784
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 10.410784129214001
$data[3,13] = 1.05358901590432
$data[3,14] = 0
$data[3,15] = 0
$data[3,16] = 0
$data[3,17] = 2.8308305572839298
$data[3,18] = 0
$data[3,19] = $true
$data[4,0] = 257
$data[4,1] = 0
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 10.1030826800597
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0
$data[4,9] = 7.9649584843580401
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 2.6020614777392099
$data[4,14] = 0
$data[4,15] = 0
$data[4,16] = 0
$data[4,17] = 2.1648755178905801
$data[4,18] = 0
$data[4,19] = $false
$data[5,0] = 258
$data[5,1] = 0
$data[5,2] = 9.3877584210417897
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 12.883147416252999
$data[5,13] = 2.06760723683978
$data[5,14] = 0
$data[5,15] = 0
$data[5,16] = 0
$data[5,17] = 1.1750880382761
$data[5,18] = 0
$data[5,19] = $false
$data[6,0] = 259
$data[6,1] = 0
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0
$data[6,9] = 12.222394850928501
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 8.5352369347823291
$data[6,13] = 1.54428685250946
$data[6,14] = 0
$data[6,15] = 0
$data[6,16] = 0
$data[6,17] = 2.9954382306213598
$data[6,18] = 0
$data[6,19] = $true
$data[7,0] = 260
$data[7,1] = 0
$data[7,2] = 0
$data[7,3] = 0
$data[7,4] = 7.2676482967029701
$data[7,5] = 8.0636969501242497
$data[7,6] = 0
$data[7,7] = 0
$data[7,8] = 0
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 2.9059863391099001
$data[7,14] = 0
$data[7,15] = 0
$data[7,16] = 0
$data[7,17] = 5.08387939765149
$data[7,18] = 0
$data[7,19] = $false
$data[8,0] = 261
$data[8,1] = 0
$data[8,2] = 0
$data[8,3] = 0
$data[8,4] = 0
$data[8,5] = 0
$data[8,6] = 9.2089215699919595
$data[8,7] = 0
$data[8,8] = 0
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 11.821111030922999
$data[8,13] = 1.5686187673729901
$data[8,14] = 0
$data[8,15] = 0
$data[8,16] = 0
$data[8,17] = 2.7607690305763799
$data[8,18] = 0
$data[8,19] = $false
$data[9,0] = 262
$data[9,1] = 0
$data[9,2] = 0
$data[9,3] = 0
$data[9,4] = 0
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 6.7817476999314001
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 9.1778771151222909
$data[9,13] = 1.9798311565906701
$data[9,14] = 0
$data[9,15] = 0
$data[9,16] = 0
$data[9,17] = 1.9722066463599399
$data[9,18] = 0
$data[9,19] = $false
$data[10,0] = 263
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 12.5558385585202
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 12.486452855583799
$data[10,13] = 1.4310801230649699
$data[10,14] = 0
$data[10,15] = 0
$data[10,16] = 0
$data[10,17] = 2.5276220355432799
$data[10,18] = 0
$data[10,19] = $false
$data[11,0] = 264
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 0
$data[11,5] = 10.8515169605094
$data[11,6] = 0
$data[11,7] = 0
$data[11,8] = 0
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 12.8603566290085
$data[11,13] = 2.6354821716309602
$data[11,14] = 0
$data[11,15] = 0
$data[11,16] = 0
$data[11,17] = 3.3049294843491701
$data[11,18] = 0
$data[11,19] = $true
$data[12,0] = 265
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 0
$data[12,4] = 8.7030684296429399
$data[12,5] = 0
$data[12,6] = 0
$data[12,7] = 0
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 8.6297365481994408
$data[12,13] = 2.5793470799381399
$data[12,14] = 0
$data[12,15] = 0
$data[12,16] = 0
$data[12,17] = 4.3999128866071002
$data[12,18] = 0
$data[12,19] = $false
$data[13,0] = 266
$data[13,1] = 0
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 13.3222133000941
$data[13,5] = 6.7216799284721196
$data[13,6] = 0
$data[13,7] = 0
$data[13,8] = 0
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 1.2297841002854599
$data[13,14] = 0
$data[13,15] = 0
$data[13,16] = 0
$data[13,17] = 4.5771904745483099
$data[13,18] = 0
$data[13,19] = $false
$data[14,0] = 267
$data[14,1] = 0
$data[14,2] = 0
$data[14,3] = 0
$data[14,4] = 0
$data[14,5] = 0
$data[14,6] = 0
$data[14,7] = 0
$data[14,8] = 0
$data[14,9] = 8.5500400063464692
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 11.251948716217299
$data[14,13] = 0.99195074394791105
$data[14,14] = 0
$data[14,15] = 0
$data[14,16] = 0
$data[14,17] = 2.4317178388836602
$data[14,18] = 0
$data[14,19] = $true
$data[15,0] = 268
$data[15,1] = 0
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 8.50924547642866
$data[15,5] = 0
$data[15,6] = 0
$data[15,7] = 0
$data[15,8] = 0
$data[15,9] = 11.530174839340701
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 2.2487668624992101
$data[15,14] = 0
$data[15,15] = 0
$data[15,16] = 0
$data[15,17] = 1.39367111609324
$data[15,18] = 0
$data[15,19] = $true
$data[16,0] = 269
$data[16,1] = 0
$data[16,2] = 0
$data[16,3] = 0
$data[16,4] = 0
$data[16,5] = 9.7566016571435199
$data[16,6] = 0
$data[16,7] = 0
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 11.4170259497181
$data[16,13] = 1.64730525746613
$data[16,14] = 0
$data[16,15] = 0
$data[16,16] = 0
$data[16,17] = 4.2178935517994001
$data[16,18] = 0
$data[16,19] = $false
$data[17,0] = 270
$data[17,1] = 0
$data[17,2] = 0
$data[17,3] = 0
$data[17,4] = 0
$data[17,5] = 10.844081307943901
$data[17,6] = 0
$data[17,7] = 0
$data[17,8] = 0
$data[17,9] = 12.5849371700888
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 1.8452484011755299
$data[17,14] = 0
$data[17,15] = 0
$data[17,16] = 0
$data[17,17] = 3.7493095003964201
$data[17,18] = 0
$data[17,19] = $false
$data[18,0] = 271
$data[18,1] = 0
$data[18,2] = 0
$data[18,3] = 0
$data[18,4] = 0
$data[18,5] = 0
$data[18,6] = 11.7073117612775
$data[18,7] = 0
$data[18,8] = 0
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 12.844216032162199
$data[18,13] = 1.0145569582829399
$data[18,14] = 0
$data[18,15] = 0
$data[18,16] = 0
$data[18,17] = 1.34607063762339
$data[18,18] = 0
$data[18,19] = $true
$data[19,0] = 272
$data[19,1] = 0
$data[19,2] = 0
$data[19,3] = 0
$data[19,4] = 0
$data[19,5] = 0
$data[19,6] = 0
$data[19,7] = 7.9663250903729201
$data[19,8] = 0
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 8.2400957005273892
$data[19,13] = 1.60740317184714
$data[19,14] = 0
$data[19,15] = 0
$data[19,16] = 0
$data[19,17] = 1.54870235671367
$data[19,18] = 0
$data[19,19] = $true
$data[20,0] = 273
$data[20,1] = 0
$data[20,2] = 0
$data[20,3] = 0
$data[20,4] = 12.067036914070799
$data[20,5] = 8.19803190383276
$data[20,6] = 0
$data[20,7] = 0
$data[20,8] = 0
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 1.9021101345115501
$data[20,14] = 0
$data[20,15] = 0
$data[20,16] = 0
$data[20,17] = 2.2800211249326399
$data[20,18] = 0
$data[20,19] = $true
$data[21,0] = 274
$data[21,1] = 0
$data[21,2] = 0
$data[21,3] = 0
$data[21,4] = 12.395699611050601
$data[21,5] = 0
$data[21,6] = 0
$data[21,7] = 0
$data[21,8] = 0
$data[21,9] = 10.549177265007501
$data[21,10] = 0
$data[21,11] = 0
$data[21,12] = 0
$data[21,13] = 2.9216195740444801
$data[21,14] = 0
$data[21,15] = 0
$data[21,16] = 0
$data[21,17] = 3.56080038613627
$data[21,18] = 0
$data[21,19] = $false
$data[22,0] = 275
$data[22,1] = 0
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 9.8631375515323398
$data[22,5] = 12.375013646398701
$data[22,6] = 0
$data[22,7] = 0
$data[22,8] = 0
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 0
$data[22,13] = 2.4277404392616901
$data[22,14] = 0
$data[22,15] = 0
$data[22,16] = 0
$data[22,17] = 1.7656294103721699
$data[22,18] = 0
$data[22,19] = $true
$data[23,0] = 276
$data[23,1] = 0
$data[23,2] = 0
$data[23,3] = 0
$data[23,4] = 0
$data[23,5] = 0
$data[23,6] = 0
$data[23,7] = 10.334694882431201
$data[23,8] = 0
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 9.2808018797427998
$data[23,13] = 2.2617059382718798
$data[23,14] = 0
$data[23,15] = 0
$data[23,16] = 0
$data[23,17] = 2.2852524658483402
$data[23,18] = 0
$data[23,19] = $false
$ws.Range("A254:T277").Value = $data

# --- 4. Columns U (Turbidity/NTU), V (Viscosity_10/cP) and W (Newtonian) follow the sheet's
#        existing pattern: a formula IF(Tn=FALSE,"NA","") unless the stability test (T) is
#        TRUE, in which case the measured values (U/V numbers, W boolean) were entered by hand.
$ws.Cells.Item(254, 21).Value = 19
$ws.Cells.Item(254, 22).Value = 5
$ws.Cells.Item(254, 23).Value = $true
$ws.Cells.Item(255, 21).Value = 21
$ws.Cells.Item(255, 22).Value = 9
$ws.Cells.Item(255, 23).Value = $true
$ws.Cells.Item(256, 21).Formula = '=IF(T256=FALSE, "NA", "")'
$ws.Cells.Item(256, 22).Formula = '=IF(T256=FALSE, "NA", "")'
$ws.Cells.Item(256, 23).Formula = '=IF(T256=FALSE, "NA", "")'
$ws.Cells.Item(257, 21).Value = 25
$ws.Cells.Item(257, 22).Value = 340
$ws.Cells.Item(257, 23).Value = $false
$ws.Cells.Item(258, 21).Formula = '=IF(T258=FALSE, "NA", "")'
$ws.Cells.Item(258, 22).Formula = '=IF(T258=FALSE, "NA", "")'
$ws.Cells.Item(258, 23).Formula = '=IF(T258=FALSE, "NA", "")'
$ws.Cells.Item(259, 21).Formula = '=IF(T259=FALSE, "NA", "")'
$ws.Cells.Item(259, 22).Formula = '=IF(T259=FALSE, "NA", "")'
$ws.Cells.Item(259, 23).Formula = '=IF(T259=FALSE, "NA", "")'
$ws.Cells.Item(260, 21).Value = 19
$ws.Cells.Item(260, 22).Value = 64
$ws.Cells.Item(260, 23).Value = $true
$ws.Cells.Item(261, 21).Formula = '=IF(T261=FALSE, "NA", "")'
$ws.Cells.Item(261, 22).Formula = '=IF(T261=FALSE, "NA", "")'
$ws.Cells.Item(261, 23).Formula = '=IF(T261=FALSE, "NA", "")'
$ws.Cells.Item(262, 21).Formula = '=IF(T262=FALSE, "NA", "")'
$ws.Cells.Item(262, 22).Formula = '=IF(T262=FALSE, "NA", "")'
$ws.Cells.Item(262, 23).Formula = '=IF(T262=FALSE, "NA", "")'
$ws.Cells.Item(263, 21).Formula = '=IF(T263=FALSE, "NA", "")'
$ws.Cells.Item(263, 22).Formula = '=IF(T263=FALSE, "NA", "")'
$ws.Cells.Item(263, 23).Formula = '=IF(T263=FALSE, "NA", "")'
$ws.Cells.Item(264, 21).Formula = '=IF(T264=FALSE, "NA", "")'
$ws.Cells.Item(264, 22).Formula = '=IF(T264=FALSE, "NA", "")'
$ws.Cells.Item(264, 23).Formula = '=IF(T264=FALSE, "NA", "")'
$ws.Cells.Item(265, 21).Value = 24
$ws.Cells.Item(265, 22).Value = 357
$ws.Cells.Item(265, 23).Value = $false
$ws.Cells.Item(266, 21).Formula = '=IF(T266=FALSE, "NA", "")'
$ws.Cells.Item(266, 22).Formula = '=IF(T266=FALSE, "NA", "")'
$ws.Cells.Item(266, 23).Formula = '=IF(T266=FALSE, "NA", "")'
$ws.Cells.Item(267, 21).Formula = '=IF(T267=FALSE, "NA", "")'
$ws.Cells.Item(267, 22).Formula = '=IF(T267=FALSE, "NA", "")'
$ws.Cells.Item(267, 23).Formula = '=IF(T267=FALSE, "NA", "")'
$ws.Cells.Item(268, 21).Value = 20
$ws.Cells.Item(268, 22).Value = 23
$ws.Cells.Item(268, 23).Value = $true
$ws.Cells.Item(269, 21).Value = 31
$ws.Cells.Item(269, 22).Value = 18
$ws.Cells.Item(269, 23).Value = $true
$ws.Cells.Item(270, 21).Formula = '=IF(T270=FALSE, "NA", "")'
$ws.Cells.Item(270, 22).Formula = '=IF(T270=FALSE, "NA", "")'
$ws.Cells.Item(270, 23).Formula = '=IF(T270=FALSE, "NA", "")'
$ws.Cells.Item(271, 21).Formula = '=IF(T271=FALSE, "NA", "")'
$ws.Cells.Item(271, 22).Formula = '=IF(T271=FALSE, "NA", "")'
$ws.Cells.Item(271, 23).Formula = '=IF(T271=FALSE, "NA", "")'
$ws.Cells.Item(272, 21).Value = 28
$ws.Cells.Item(272, 22).Value = 1030
$ws.Cells.Item(272, 23).Value = $false
$ws.Cells.Item(273, 21).Value = 28
$ws.Cells.Item(273, 22).Value = 170
$ws.Cells.Item(273, 23).Value = $false
$ws.Cells.Item(274, 21).Value = 25
$ws.Cells.Item(274, 22).Value = 60
$ws.Cells.Item(274, 23).Value = $true
$ws.Cells.Item(275, 21).Formula = '=IF(T275=FALSE, "NA", "")'
$ws.Cells.Item(275, 22).Formula = '=IF(T275=FALSE, "NA", "")'
$ws.Cells.Item(275, 23).Formula = '=IF(T275=FALSE, "NA", "")'
$ws.Cells.Item(276, 21).Value = 29
$ws.Cells.Item(276, 22).Value = 15
$ws.Cells.Item(276, 23).Value = $true
$ws.Cells.Item(277, 21).Formula = '=IF(T277=FALSE, "NA", "")'
$ws.Cells.Item(277, 22).Formula = '=IF(T277=FALSE, "NA", "")'
$ws.Cells.Item(277, 23).Formula = '=IF(T277=FALSE, "NA", "")'

# --- 5. Recalculate so the new formula cells carry a cached value, then refresh the used range,
#        the conditional-formatting range and the on-screen selection to match the new extent.
$excel.CalculateFullRebuild()

$fc = $ws.Range("A1:W253").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:W277"))

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("A1:W277").Select()

# --- 6. Old_ProtocolNotes sheet: cursor moved to A4 as part of the same editing session.
$notes = $wb.Worksheets.Item("Old_ProtocolNotes")
$notes.Activate()
$notes.Range("A4").Select()
$ws.Activate()

